$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the diff: (cellRef, newValue, forceText)
$updates = @(
    ,@('D2', '64.060.15', $false)
    ,@('E2', '  -3.61%  ', $false)
    ,@('D3', '3.160.71', $false)
    ,@('E3', '  -8.83%  ', $false)
    ,@('E4', '  +0.00%  ', $false)
    ,@('D5', '562.04', $true)
    ,@('E5', '  -4.26%  ', $false)
    ,@('D6', '168.32', $true)
    ,@('E6', '  -4.85%  ', $false)
    ,@('E7', '  +0.04%  ', $false)
    ,@('D8', '0.601', $true)
    ,@('E8', '  -3.42%  ', $false)
    ,@('D9', '3.159.61', $false)
    ,@('E9', '  -8.87%  ', $false)
    ,@('D10', '0.124', $true)
    ,@('E10', '  -7.39%  ', $false)
    ,@('D11', '6.60', $true)
    ,@('E11', '  -5.20%  ', $false)
    ,@('D12', '0.394', $true)
    ,@('E12', '  -5.40%  ', $false)
    ,@('D13', '3.709.95', $false)
    ,@('E13', '  -8.81%  ', $false)
    ,@('E14', '  +1.50%  ', $false)
    ,@('D15', '27.26', $true)
    ,@('E15', '  -10.21%  ', $false)
    ,@('D16', '64.044.97', $false)
    ,@('E16', '  -3.47%  ', $false)
    ,@('D17', '0.0000163', $true)
    ,@('E17', '  -5.87%  ', $false)
    ,@('D18', '3.169.58', $false)
    ,@('E18', '  -8.43%  ', $false)
    ,@('D19', '5.73', $true)
    ,@('E19', '  -4.20%  ', $false)
    ,@('D20', '12.92', $true)
    ,@('E20', '  -6.61%  ', $false)
    ,@('D21', '351.09', $true)
    ,@('E21', '  -5.94%  ', $false)
    ,@('D22', '7.15', $true)
    ,@('E22', '  -6.59%  ', $false)
    ,@('D23', '1.00', $true)
    ,@('E23', '  -0.03%  ', $false)
    ,@('D24', '68.61', $true)
    ,@('E24', '  -6.37%  ', $false)
    ,@('D25', '0.0000118', $true)
    ,@('E25', '  -6.98%  ', $false)
    ,@('D26', '0.501', $true)
    ,@('E26', '  -6.28%  ', $false)
    ,@('D27', '9.48', $true)
    ,@('E27', '  -5.01%  ', $false)
    ,@('E28', '  -1.30%  ', $false)
    ,@('E29', '  +0.30%  ', $false)
    ,@('D30', '0.998', $true)
    ,@('E30', '  -0.20%  ', $false)
    ,@('D31', '5.50', $true)
    ,@('E31', '  -7.01%  ', $false)
    ,@('D32', '1.89', $true)
    ,@('E32', '  -5.68%  ', $false)
    ,@('D33', '21.89', $true)
    ,@('E33', '  -7.70%  ', $false)
    ,@('B34', 'Aptos', $false)
    ,@('C34', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', $false)
    ,@('D34', '6.61', $true)
    ,@('E34', '  -6.38%  ', $false)
    ,@('B35', 'Fetch.AI', $false)
    ,@('C35', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', $false)
    ,@('D35', '1.20', $true)
    ,@('E35', '  -5.90%  ', $false)
    ,@('D36', '1.43', $true)
    ,@('E36', '  -8.01%  ', $false)
    ,@('D37', '153.57', $true)
    ,@('E37', '  -4.67%  ', $false)
    ,@('D38', '0.813', $true)
    ,@('E38', '  -8.38%  ', $false)
    ,@('D39', '25.69', $true)
    ,@('E39', '  -9.31%  ', $false)
    ,@('D40', '2.51', $true)
    ,@('E40', '  -3.67%  ', $false)
    ,@('D41', '1.69', $true)
    ,@('E41', '  -7.14%  ', $false)
    ,@('D42', '2.595.64', $false)
    ,@('E42', '  -6.71%  ', $false)
    ,@('D43', '4.16', $true)
    ,@('E43', '  -7.71%  ', $false)
    ,@('B44', 'OKB', $false)
    ,@('C44', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', $false)
    ,@('D44', '39.35', $true)
    ,@('E44', '  -1.70%  ', $false)
    ,@('B45', 'RenderToken', $false)
    ,@('C45', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', $false)
    ,@('D45', '5.93', $true)
    ,@('E45', '  -8.34%  ', $false)
    ,@('D46', '0.0650', $true)
    ,@('E46', '  -6.49%  ', $false)
    ,@('D47', '23.57', $true)
    ,@('E47', '  -7.33%  ', $false)
    ,@('D48', '317.46', $true)
    ,@('E48', '  -6.62%  ', $false)
    ,@('D49', '0.0268', $true)
    ,@('E49', '  -8.61%  ', $false)
    ,@('E50', '  -3.30%  ', $false)
    ,@('E51', '  -0.07%  ', $false)
)

foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    $forceText = $u[2]
    $cell = $ws.Range($ref)
    if ($forceText) {
        # Prevent Excel from re-interpreting a plain-number-looking string
        # as a numeric value; keep it stored as text exactly as scraped,
        # then restore the default (Normal) cell style so no formatting
        # side effects are introduced.
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}

